$d = $word.ActiveDocument

# --- Edit 1 -----------------------------------------------------------
# Locate the long "linear regression hyper-parameters" paragraph via
# Find, then resolve which document paragraph contains the match and
# move to the paragraph right after it. That paragraph is currently
# empty (a single bare run holding just <w:rtl w:val="0"/>). Give it
# the body text size (12pt / 24 half-points for both sz and szCs) and
# fill in the missing R2-evaluation sentence.
$findRange = $d.Content
$found = $findRange.Find.Execute(
    "За степен учења је првобитно одабрана вредност",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$matchStart = $findRange.Start

$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Start -le $matchStart -and $matchStart -lt $cand.Range.End) {
        $anchorIndex = $i
        break
    }
}

$targetPara = $d.Paragraphs.Item($anchorIndex + 1)

$targetRange = $targetPara.Range
$targetRange.Text = "За евалуацију модела коришћена је R2 метрика. За конкретну имплементацију, најбољи резултат је износио 0.61."

$targetRange2 = $targetPara.Range
$targetRange2.Font.Size = 12
$targetRange2.Font.SizeBi = 12

# --- Edit 2 -------------------------------------------------------------
# Append a brand-new paragraph at the very end of the document body
# (right before sectPr). It inherits the formatting of the paragraph
# that currently ends the document (ind left/firstLine 0, sz/szCs 24),
# and gets the new "precision" evaluation sentence.
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Range.Text = "За евалуацију модела коришћена је прецизност која је за ову имплементацију износила око 65%."
